$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.410.03"
$ws.Range("E2").Value = "'  +1.07%  "
$ws.Range("D3").Value = "'1.672.88"
$ws.Range("E3").Value = "'  +1.07%  "
$ws.Range("E4").Value = "'  +0.86%  "
$ws.Range("D5").Value = "'220.89"
$ws.Range("E5").Value = "'  +1.41%  "
$ws.Range("D6").Value = "'0.5364"
$ws.Range("E6").Value = "'  +1.04%  "
$ws.Range("E8").Value = "'  +2.43%  "
$ws.Range("D9").Value = "'0.06419"
$ws.Range("E9").Value = "'  +1.44%  "
$ws.Range("D10").Value = "'21.04"
$ws.Range("E10").Value = "'  +3.00%  "
$ws.Range("D11").Value = "'0.07851"
$ws.Range("E11").Value = "'  +0.79%  "
$ws.Range("D12").Value = "'4.570"
$ws.Range("E12").Value = "'  +1.04%  "
$ws.Range("D13").Value = "'1.673.20"
$ws.Range("E13").Value = "'  -0.58%  "
$ws.Range("D14").Value = "'1.903.01"
$ws.Range("E14").Value = "'  +1.08%  "
$ws.Range("D15").Value = "'0.5651"
$ws.Range("E15").Value = "'  +2.95%  "
$ws.Range("D16").Value = "'0.0₅8206"
$ws.Range("D17").Value = "'66.52"
$ws.Range("E17").Value = "'  +1.67%  "
$ws.Range("D18").Value = "'26.449.60"
$ws.Range("E18").Value = "'  +1.26%  "
$ws.Range("E19").Value = "'  +0.71%  "
$ws.Range("D20").Value = "'4.722"
$ws.Range("E20").Value = "'  +2.61%  "
$ws.Range("D21").Value = "'197.98"
$ws.Range("E21").Value = "'  +3.64%  "
$ws.Range("E22").Value = "'  +2.81%  "
$ws.Range("D23").Value = "'6.081"
$ws.Range("E23").Value = "'  +0.93%  "
$ws.Range("D25").Value = "'146.58"
$ws.Range("E25").Value = "'  +0.79%  "
$ws.Range("D26").Value = "'0.1235"
$ws.Range("E26").Value = "'  +0.41%  "
$ws.Range("E27").Value = "'  +0.80%  "
$ws.Range("D28").Value = "'16.28"
$ws.Range("E28").Value = "'  +1.92%  "
$ws.Range("D29").Value = "'1.515"
$ws.Range("D30").Value = "'0.05897"
$ws.Range("E30").Value = "'  +1.97%  "
$ws.Range("D31").Value = "'1.291"
$ws.Range("E31").Value = "'  +1.38%  "
$ws.Range("D32").Value = "'3.585"
$ws.Range("E32").Value = "'  +0.75%  "
$ws.Range("D33").Value = "'3.313"
$ws.Range("E33").Value = "'  +1.16%  "
$ws.Range("D34").Value = "'1.626"
$ws.Range("E34").Value = "'  +1.35%  "
$ws.Range("D35").Value = "'0.9725"
$ws.Range("E35").Value = "'  +2.21%  "
$ws.Range("E36").Value = "'  +1.91%  "
$ws.Range("D37").Value = "'2.432"
$ws.Range("E37").Value = "'  +0.68%  "
$ws.Range("D38").Value = "'0.5843"
$ws.Range("E38").Value = "'  +1.56%  "
$ws.Range("E39").Value = "'  +0.25%  "
$ws.Range("D40").Value = "'1.079.81"
$ws.Range("E40").Value = "'  +3.92%  "
$ws.Range("D41").Value = "'5.924"
$ws.Range("E41").Value = "'  +2.24%  "
$ws.Range("D42").Value = "'0.8669"
$ws.Range("E42").Value = "'  +1.22%  "
$ws.Range("D44").Value = "'104.47"
$ws.Range("E44").Value = "'  -0.24%  "
$ws.Range("D45").Value = "'1.812.13"
$ws.Range("E45").Value = "'  +0.91%  "
$ws.Range("E46").Value = "'  +2.59%  "
$ws.Range("E47").Value = "'  -4.04%  "
$ws.Range("E48").Value = "'  +0.98%  "
$ws.Range("D49").Value = "'0.4407"
$ws.Range("E49").Value = "'  +1.68%  "
$ws.Range("D50").Value = "'8.091"
$ws.Range("E50").Value = "'  +2.89%  "
$ws.Range("E51").Value = "'  +0.35%  "
